# Apply the edits described by the diff:
# 1. Rename the sheet "C_11" -> "C_11.2"
# 2. Apply a date number format (built-in numFmtId 15, "d-mmm-yy") to cell B2
#    (keeping the rest of its formatting, i.e. the bold font, unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "C_11.2"

# Apply date number format to B2
$ws.Range("B2").NumberFormat = "d-mmm-yy"
